$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet tab
$ws.Name = "Jayant Yadav"

# Insert a new column before column A, shifting existing data (teamName..result) to B..M
$ws.Columns.Item(1).Insert()

# New header for the inserted column
$ws.Range("A1").Value = "matchNo"

# New matchNo value for the existing data row
$ws.Range("A2").Value = "46th"

# Append two new rows of scraped data
$ws.Range("A3").Value = "13th"
$ws.Range("B3").Value = "Mumbai Indians"
$ws.Range("C3").Value = "Jayant Yadav"
$ws.Range("D3").Value = "c & b Rabada"
$ws.Range("E3").Value = "23"
$ws.Range("F3").Value = "22"
$ws.Range("G3").Value = "1"
$ws.Range("H3").Value = "0"
$ws.Range("I3").Value = "104.54"
$ws.Range("J3").Value = "Delhi Capitals"
$ws.Range("K3").Value = "Chennai"
$ws.Range("L3").Value = "April 20"
$ws.Range("M3").Value = "Capitals won by 6 wickets (with 5 balls remaining)"

$ws.Range("A4").Value = "17th"
$ws.Range("B4").Value = "Mumbai Indians"
$ws.Range("C4").Value = "Jayant Yadav"
$ws.Range("D4").Value = ""
$ws.Range("E4").Value = "0"
$ws.Range("F4").Value = "0"
$ws.Range("G4").Value = "0"
$ws.Range("H4").Value = "0"
$ws.Range("I4").Value = "-"
$ws.Range("J4").Value = "Punjab Kings"
$ws.Range("K4").Value = "Chennai"
$ws.Range("L4").Value = "April 23"
$ws.Range("M4").Value = "Punjab Kings won by 9 wickets (with 14 balls remaining)"
